# Reorders the comma-separated "Recorded By" entries in column G.
# Rule: any element that is exactly the lowercase literal "system" stays
# fixed in its position; all other elements in the list have their
# relative order reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy($s) {
    $parts = $s -split ", "
    $n = $parts.Count

    if ($n -le 1) {
        return $s
    }

    $restIdx = @()
    for ($i = 0; $i -lt $n; $i++) {
        if (-not $parts[$i].Equals("system")) {
            $restIdx += $i
        }
    }

    $restVals = @()
    foreach ($i in $restIdx) {
        $restVals += $parts[$i]
    }

    $m = $restVals.Count
    $reversed = @()
    for ($i = $m - 1; $i -ge 0; $i--) {
        $reversed += $restVals[$i]
    }

    for ($k = 0; $k -lt $restIdx.Count; $k++) {
        $parts[$restIdx[$k]] = $reversed[$k]
    }

    return ($parts -join ", ")
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $newVal = Transform-RecordedBy $val
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
